# Apply coverage value updates to the "Platform Coverage" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Make this the active sheet/window so the selection/view changes stick.
$ws.Activate()

# Update existing coverage values in row 2 from 0.75 to 0.6
$ws.Range("H2").Value = 0.6
$ws.Range("J2").Value = 0.6
$ws.Range("L2").Value = 0.6
$ws.Range("N2").Value = 0.6
$ws.Range("P2").Value = 0.6

# Add new coverage values (0.5) across rows 3 and 4
$cols = @("R", "T", "V", "X", "Z", "AB", "AD")
foreach ($col in $cols) {
    $ws.Range($col + "3").Value = 0.5
    $ws.Range($col + "4").Value = 0.5
}

# Update the view: scroll back to A1 (remove frozen/scrolled topLeftCell at L1)
# and move the selection to AD3:AD4, matching active cell AD3.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("AD3:AD4").Select()
